# Updated cryptos list on Sat May 18 02:55:23 UTC 2024 with GitHub Actions
# Refresh price/volume columns (and a couple of row re-orderings) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.901.65'
$ws.Range('E2').Value = '  +2.29%  '
$ws.Range('D3').Value = '3.098.41'
$ws.Range('E3').Value = '  +4.97%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'579.46"
$ws.Range('E5').Value = '  +1.54%  '
$ws.Range('D6').Value = "'172.97"
$ws.Range('E6').Value = '  +7.66%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.094.73'
$ws.Range('E8').Value = '  +4.95%  '
$ws.Range('E9').Value = '  +1.26%  '
$ws.Range('E10').Value = '  -2.60%  '
$ws.Range('E11').Value = '  +4.05%  '
$ws.Range('D12').Value = "'0.481"
$ws.Range('E12').Value = '  +4.80%  '
$ws.Range('E13').Value = '  +2.31%  '
$ws.Range('D14').Value = "'37.10"
$ws.Range('E14').Value = '  +7.26%  '
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('D16').Value = '3.609.76'
$ws.Range('E16').Value = '  +4.92%  '
$ws.Range('D17').Value = '66.863.51'
$ws.Range('E17').Value = '  +2.25%  '
$ws.Range('E18').Value = '  +2.32%  '
$ws.Range('D19').Value = '3.097.34'
$ws.Range('E19').Value = '  +4.98%  '
$ws.Range('D20').Value = "'16.22"
$ws.Range('E20').Value = '  +1.78%  '
$ws.Range('D21').Value = "'481.80"
$ws.Range('E21').Value = '  +8.15%  '
$ws.Range('E22').Value = '  +2.56%  '
$ws.Range('E23').Value = '  +2.93%  '
$ws.Range('D24').Value = "'84.04"
$ws.Range('E24').Value = '  +2.04%  '
$ws.Range('E25').Value = '  +4.71%  '
$ws.Range('D26').Value = "'13.00"
$ws.Range('E26').Value = '  +6.37%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = "'10.01"
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').Value = "'1.00"
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').Value = "'2.38"
$ws.Range('E30').Value = '  -3.42%  '
$ws.Range('E31').Value = '  +3.33%  '
$ws.Range('E32').Value = '  -0.65%  '
$ws.Range('D33').Value = "'28.72"
$ws.Range('E33').Value = '  +5.81%  '
$ws.Range('E34').Value = '  +1.89%  '
$ws.Range('D35').Value = "'1.00"
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  +3.21%  '
$ws.Range('E37').Value = '  +2.46%  '
$ws.Range('D38').Value = "'47.95"
$ws.Range('E38').Value = '  +6.85%  '
$ws.Range('E39').Value = '  +8.01%  '
$ws.Range('E40').Value = '  +2.19%  '
$ws.Range('E41').Value = '  +5.05%  '
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('D43').Value = "'8.65"
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('E44').Value = '  -1.85%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.818.14'
$ws.Range('E45').Value = '  +5.06%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = "'0.0359"
$ws.Range('E46').Value = '  +2.15%  '
$ws.Range('D47').Value = "'378.79"
$ws.Range('E47').Value = '  -1.10%  '
$ws.Range('D48').Value = "'135.26"
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').Value = "'24.78"
$ws.Range('E50').Value = '  +4.59%  '
$ws.Range('E51').Value = '  +1.90%  '

$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
